$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.254.90"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "3.678.32"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "681.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("D7").Value = "3.676.52"
$ws.Range("E7").Value = "  -3.61%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.90%  "
$ws.Range("E10").Value = "  -8.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.61%  "
$ws.Range("D15").Value = "4.299.87"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "3.676.23"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").Value = "69.272.91"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  -6.39%  "
$ws.Range("E20").Value = "  -7.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "480.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.07%  "
$ws.Range("E22").Value = "  -7.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.665"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.11%  "
$ws.Range("D25").Value = "3.825.21"
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000130"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.29%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.15%  "
$ws.Range("E30").Value = "  -10.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("E35").Value = "  -7.05%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "3.648.48"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0933"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.14%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "162.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("E47").Value = "  -13.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000287"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.20%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -4.50%  "
